$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '28.445.34'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.22%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.552.47'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.94%  '

$ws.Range("E4").Value = '  -0.35%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '210.78'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.484'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.65%  '

$ws.Range("E7").Value = '  -0.29%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '24.17'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.44%  '

$ws.Range("E9").Value = '  -1.69%  '

$ws.Range("E10").Value = '  -1.08%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0892'
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.774.74'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.86%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.545.12'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.32%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '28.445.16'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.09%  '

$ws.Range("E15").Value = '  -2.04%  '

$ws.Range("E16").Value = '  -1.50%  '

$ws.Range("E17").Value = '  -1.40%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '229.02'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("E19").Value = '  -1.22%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0₃0672'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.36%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.26%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '3.90'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.60%  '

$ws.Range("E23").Value = '  -2.47%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.06'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.21%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '150.95'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.26%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '14.73'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.95%  '

$ws.Range("E27").Value = '  -1.35%  '

$ws.Range("E28").Value = '  -0.28%  '

$ws.Range("E29").Value = '  -2.78%  '

$ws.Range("E30").Value = '  -3.88%  '

$ws.Range("E31").Value = '  -4.46%  '

$ws.Range("E32").Value = '  -1.30%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.388.84'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.89%  '

$ws.Range("E34").Value = '  -3.00%  '

$ws.Range("E35").Value = '  -5.16%  '

$ws.Range("E36").Value = '  -1.74%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.29'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.88%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.65'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.42%  '

$ws.Range("E39").Value = '  -1.55%  '

$ws.Range("E40").Value = '  +3.87%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.515'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.09%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("E43").Value = '  -2.16%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0464'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.94%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '5.34'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.76%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '61.98'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.21%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.687.88'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.81%  '

$ws.Range("E48").Value = '  -6.30%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '85.40'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.14%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '43.27'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +7.17%  '

$ws.Range("E51").Value = '  +1.14%  '
